# Applies the "slice and dice" data-fill edit:
#  - Adds/overwrites the jan_2015..jun_2015 columns (O:T) for rows 2-13
#  - Best-effort update of the saved window size recorded in the workbook view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: record the window geometry the way Excel would when the
# document window is resized/maximized before saving.
$win = $excel.ActiveWindow
$win.Width = 20490
$win.Height = 7755

# New values for columns O (jan_2015) through T (jun_2015), rows 2-13.
$newData = @{
    2  = @(225, 274, 260, 249, 248, 126)
    3  = @(225, 260, 248, 253, 250, 130)
    4  = @(234, 255, 265, 267, 254, 135)
    5  = @(244, 258, 254, 243, 244, 125)
    6  = @(242, 265, 241, 274, 245, 137)
    7  = @(236, 272, 256, 261, 258, 139)
    8  = @(229, 255, 250, 260, 268, 136)
    9  = @(229, 267, 236, 279, 249, 129)
    10 = @(226, 258, 258, 257, 260, 124)
    11 = @(229, 255, 255, 258, 232, 129)
    12 = @(245, 279, 251, 261, 233, 127)
    13 = @(231, 263, 265, 243, 264, 127)
}

$columns = @("O", "P", "Q", "R", "S", "T")

foreach ($row in $newData.Keys) {
    $values = $newData[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$row"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
